$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab / sheet name to reflect the new "through" date
$ws.Name = "Through 2022-07-26"

# Update the header label for the 2022 column (row 1, column I)
$ws.Range("I1").Value = "2022 (through 07-26)"

# Update July's 2022 value (row 8, column I)
$ws.Range("I8").Value = 146

# Update the Total 2022 value (row 14, column I)
$ws.Range("I14").Value = 952
